# Add a "TargetType" column to the Skill sheet of SkillData.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill")

# Insert a new column before the existing ATKMul column (old column C)
$ws.Columns("C").Insert()

# New column width (closest reproducible value to the authored 11.125 width)
$ws.Columns("C").ColumnWidth = 10.428571428571429

# Header
$ws.Range("C1").Value = "TargetType"

# Data values for the new TargetType column
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("C17").Value = 2
$ws.Range("C18").Value = 2
$ws.Range("C19").Value = 2
$ws.Range("C20").Value = 2
$ws.Range("C21").Value = 2
$ws.Range("C22").Value = 2

# Match the author's final selection/view state
$ws.Activate()
$ws.Range("C1").Select()
